# Updates cryptos list price/volume columns (D/E) for rows 2-51
# to match the refreshed data in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.792.83"
$ws.Range("E2").Value = "  -3.94%  "
$ws.Range("D3").Value = "1.816.41"
$ws.Range("E3").Value = "  -3.09%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'276.63"
$ws.Range("E5").Value = "  -8.35%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").Value = "'0.5116"
$ws.Range("E7").Value = "  -4.90%  "
$ws.Range("D8").Value = "'0.3516"
$ws.Range("E8").Value = "  -6.48%  "
$ws.Range("D9").Value = "'44.83"
$ws.Range("E9").Value = "  -1.77%  "
$ws.Range("D10").Value = "'0.06666"
$ws.Range("E10").Value = "  -7.29%  "
$ws.Range("D11").Value = "'20.03"
$ws.Range("E11").Value = "  -7.21%  "
$ws.Range("D12").Value = "'0.8305"
$ws.Range("E12").Value = "  -6.55%  "
$ws.Range("D13").Value = "'0.07857"
$ws.Range("E13").Value = "  -3.82%  "
$ws.Range("D14").Value = "1.814.48"
$ws.Range("E14").Value = "  -3.49%  "
$ws.Range("D15").Value = "'5.079"
$ws.Range("E15").Value = "  -3.64%  "
$ws.Range("D16").Value = "'87.51"
$ws.Range("E16").Value = "  -6.24%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("D18").Value = "'14.12"
$ws.Range("E18").Value = "  -4.31%  "
$ws.Range("D19").Value = "'0.000008023"
$ws.Range("E19").Value = "  -6.25%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").Value = "25.866.43"
$ws.Range("E21").Value = "  -3.82%  "
$ws.Range("D22").Value = "'4.724"
$ws.Range("E22").Value = "  -5.22%  "
$ws.Range("E23").Value = "  -6.50%  "
$ws.Range("D24").Value = "'6.084"
$ws.Range("E24").Value = "  -4.83%  "
$ws.Range("D25").Value = "'141.25"
$ws.Range("E25").Value = "  -4.12%  "
$ws.Range("D26").Value = "'2.191"
$ws.Range("E26").Value = "  -3.31%  "
$ws.Range("D27").Value = "'1.674"
$ws.Range("E27").Value = "  -3.75%  "
$ws.Range("D28").Value = "'17.08"
$ws.Range("E28").Value = "  -5.33%  "
$ws.Range("D29").Value = "'109.42"
$ws.Range("E29").Value = "  -4.29%  "
$ws.Range("D30").Value = "'4.351"
$ws.Range("E30").Value = "  -8.14%  "
$ws.Range("D31").Value = "'4.236"
$ws.Range("E31").Value = "  -7.91%  "
$ws.Range("D32").Value = "'0.08811"
$ws.Range("E32").Value = "  -3.85%  "
$ws.Range("E33").Value = "  -2.35%  "
$ws.Range("D34").Value = "'0.7309"
$ws.Range("E34").Value = "  -9.62%  "
$ws.Range("D35").Value = "'1.137"
$ws.Range("E35").Value = "  -3.35%  "
$ws.Range("D36").Value = "'2.875"
$ws.Range("E36").Value = "  -3.51%  "
$ws.Range("D37").Value = "'3.152"
$ws.Range("E37").Value = "  -1.11%  "
$ws.Range("D38").Value = "'1.000"
$ws.Range("E38").Value = "  +0.32%  "
$ws.Range("D39").Value = "'2.364"
$ws.Range("E39").Value = "  -9.32%  "
$ws.Range("D40").Value = "'0.5225"
$ws.Range("E40").Value = "  -14.04%  "
$ws.Range("D41").Value = "'0.01847"
$ws.Range("E41").Value = "  -5.71%  "
$ws.Range("D42").Value = "'0.9559"
$ws.Range("E42").Value = "  -11.02%  "
$ws.Range("D43").Value = "'111.72"
$ws.Range("E43").Value = "  -3.71%  "
$ws.Range("D44").Value = "'6.194"
$ws.Range("E44").Value = "  -5.94%  "
$ws.Range("D45").Value = "'8.034"
$ws.Range("E45").Value = "  -9.13%  "
$ws.Range("D46").Value = "'1.001"
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("D47").Value = "'0.4567"
$ws.Range("E47").Value = "  -11.97%  "
$ws.Range("D48").Value = "'0.1361"
$ws.Range("D49").Value = "'36.78"
$ws.Range("E49").Value = "  -2.40%  "
$ws.Range("D50").Value = "'9.227"
$ws.Range("E50").Value = "  -6.89%  "
$ws.Range("D51").Value = "'1.502"
$ws.Range("E51").Value = "  -8.14%  "
